$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.193.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "'1.824.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'234.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "'0.5989"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.06945"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.57%  "
$ws.Range("D9").Value = "'0.2762"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").Value = "'23.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.43%  "
$ws.Range("D11").Value = "'0.07598"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'1.831.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.97%  "
$ws.Range("D13").Value = "'4.728"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.83%  "
$ws.Range("D14").Value = "'0.6277"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.90%  "
$ws.Range("D15").Value = "'0.000009826"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "'77.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.35%  "
$ws.Range("D17").Value = "'29.014.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "'5.526"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.47%  "
$ws.Range("D19").Value = "'215.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.40%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").Value = "'11.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("D22").Value = "'6.848"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.35%  "
$ws.Range("D23").Value = "'1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "'155.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "'7.944"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").Value = "'16.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("D28").Value = "'0.06488"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.56%  "
$ws.Range("D29").Value = "'1.423"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("D31").Value = "'3.815"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.16%  "
$ws.Range("D32").Value = "'3.778"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("E33").Value = "  -3.43%  "
$ws.Range("D34").Value = "'1.719"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").Value = "'0.6454"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("D36").Value = "'2.542"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").Value = "'2.745"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").Value = "'0.01750"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("D39").Value = "'6.591"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("D40").Value = "'1.134.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.70%  "
$ws.Range("D41").Value = "'0.8921"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.46%  "
$ws.Range("D42").Value = "'1.003"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "'1.989.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("D44").Value = "'100.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("D45").Value = "'62.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("E46").Value = "  -2.36%  "
$ws.Range("D47").Value = "'1.615"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").Value = "'8.469"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").Value = "'0.05494"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").Value = "'0.4530"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("D51").Value = "'6.355"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.09%  "
